$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 4.592460999999999
$ws.Range("H2").Value = 13.777383
$ws.Range("I2").Value = 0.003302946473568516
$ws.Range("J2").Value = 0.003302946473568516
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 2.421679695496333
$ws.Range("R2").Value = 21.79511725946699
$ws.Range("S2").Value = 0.003302946473568516
$ws.Range("T2").Value = 0.003302946473568516

# Row 3
$ws.Range("I3").Value = 0.8667225374846176
$ws.Range("J3").Value = 0.8667225374846176
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 635.4702952203376
$ws.Range("R3").Value = 5719.232656983037
$ws.Range("S3").Value = 0.8667225374846176
$ws.Range("T3").Value = 0.8667225374846176

# Row 4
$ws.Range("G4").Value = 180.7183073333333
$ws.Range("H4").Value = 542.1549220000001
$ws.Range("I4").Value = 0.1299745160418139
$ws.Range("J4").Value = 0.1299745160418139
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 95.29571518921978
$ws.Range("R4").Value = 857.661436702978
$ws.Range("S4").Value = 0.1299745160418139
$ws.Range("T4").Value = 0.1299745160418139
